$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.28018141678364827
$ws.Range("A2").Value = -0.0059999999752911037
$ws.Range("A3").Value = -0.0039999999768713934
$ws.Range("A4").Value = -0.0079999999582334169
$ws.Range("A5").Value = -0.0029999999744223516
$ws.Range("A6").Value = -0.0019999999706250549
$ws.Range("A7").Value = -0.0099999999391879868
$ws.Range("A8").Value = -0.015458195843259581
$ws.Range("A9").Value = -0.0019999999694833015
$ws.Range("A10").Value = -0.0019999999696036497
$ws.Range("A11").Value = -0.0029999999657217558
$ws.Range("A12").Value = 0.028317199012787952
$ws.Range("A13").Value = -0.0034999999615905608
$ws.Range("A14").Value = -0.0079999999434399172
$ws.Range("A15").Value = -0.00099999997020283615
$ws.Range("A16").Value = -0.0019999999659057188
$ws.Range("A17").Value = -0.001999999965319077
$ws.Range("A18").Value = -0.0039999999573980816
$ws.Range("A19").Value = -0.0039999999822626364
$ws.Range("A20").Value = -0.0039999999743987047
$ws.Range("A21").Value = -0.069390397512012036
$ws.Range("A22").Value = -0.0039999999722963864
$ws.Range("A23").Value = -0.0049999999717131871
$ws.Range("A24").Value = -0.019999999911212818
$ws.Range("A25").Value = -0.019999999910082167
$ws.Range("A26").Value = -0.0024999999674122364
$ws.Range("A27").Value = -0.0024999999655230809
$ws.Range("A28").Value = -0.0019999999600104346
$ws.Range("A29").Value = -0.0069999999349983355
$ws.Range("A30").Value = -0.059999999726336739
$ws.Range("A31").Value = -0.0069999999302261529
$ws.Range("A32").Value = -0.0099999999182749377
$ws.Range("A33").Value = -0.0039999999414721543
